# Add new "STanH" worksheet at the end, cloned from the last sheet (FTIC2024) so it
# inherits column widths / number formats / styles, then overwrite its contents.
$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("FTIC2024")
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "STanH"

# Clear any leftover formulas/values beyond what we are about to (re)write so the sheet
# only contains exactly the 8 data rows (header + 7 rows) from the target workbook.
$ws.Range("A2:M7").ClearContents()

$ws.Range("A2").Value = "Derivation"
$ws.Range("B2").Value = 0.260041666666666
$ws.Range("C2").Value = 0.262208333333333
$ws.Range("D2").Formula = "=ABS(B2-C2)"
$ws.Range("E2").Value = 30.4812666666666
$ws.Range("F2").Value = 30.4813875
$ws.Range("G2").Formula = "=ABS(E2-F2)"
$ws.Range("H2").Value = 0.9578318870332475
$ws.Range("I2").Value = 0.95783868318231
$ws.Range("J2").Formula = "=ABS(H2-I2)"
$ws.Range("K2").Value = 56.4115757037601
$ws.Range("L2").Value = 0.125690927108128
$ws.Range("M2").Value = 0.292681276798248

$ws.Range("A3").Value = "Derivation"
$ws.Range("B3").Value = 0.315125
$ws.Range("C3").Value = 0.314791666666666
$ws.Range("D3").Formula = "=ABS(B3-C3)"
$ws.Range("E3").Value = 31.6956458333333
$ws.Range("F3").Value = 31.6955916666666
$ws.Range("G3").Formula = "=ABS(E3-F3)"
$ws.Range("H3").Value = 0.9680025464953931
$ws.Range("I3").Value = 0.9680035595361738
$ws.Range("J3").Formula = "=ABS(H3-I3)"
$ws.Range("K3").Value = 48.491065844448
$ws.Range("L3").Value = 0.10248064994812
$ws.Range("M3").Value = 0.25651100029548

$ws.Range("A4").Value = "Derivation"
$ws.Range("B4").Value = 0.360333333333333
$ws.Range("C4").Value = 0.358458333333333
$ws.Range("D4").Formula = "=ABS(B4-C4)"
$ws.Range("E4").Value = 32.4197958333333
$ws.Range("F4").Value = 32.4197874999999
$ws.Range("G4").Formula = "=ABS(E4-F4)"
$ws.Range("H4").Value = 0.9728039759431781
$ws.Range("I4").Value = 0.9728018363027288
$ws.Range("J4").Formula = "=ABS(H4-I4)"
$ws.Range("K4").Value = 43.1704180486746
$ws.Range("L4").Value = 0.0901614824930826
$ws.Range("M4").Value = 0.235784205297629

$ws.Range("A5").Value = "Anchor"
$ws.Range("B5").Value = 0.435666666666666
$ws.Range("C5").Value = 0.431791666666666
$ws.Range("D5").Formula = "=ABS(B5-C5)"
$ws.Range("E5").Value = 33.3445958333333
$ws.Range("F5").Value = 33.3449291666666
$ws.Range("G5").Formula = "=ABS(E5-F5)"
$ws.Range("H5").Value = 0.9779623854201963
$ws.Range("I5").Value = 0.9779695729079207
$ws.Range("J5").Formula = "=ABS(H5-I5)"
$ws.Range("K5").Value = 36.8194342912204
$ws.Range("L5").Value = 0.0752207338809967
$ws.Range("M5").Value = 0.214156863590081

$ws.Range("A6").Value = "Derivation"
$ws.Range("B6").Value = 0.511708333333333
$ws.Range("C6").Value = 0.506791666666666
$ws.Range("D6").Formula = "=ABS(B6-C6)"
$ws.Range("E6").Value = 33.9967208333333
$ws.Range("F6").Value = 33.9966
$ws.Range("G6").Formula = "=ABS(E6-F6)"
$ws.Range("H6").Value = 0.9810897944309005
$ws.Range("I6").Value = 0.9810883611073168
$ws.Range("J6").Formula = "=ABS(H6-I6)"
$ws.Range("K6").Value = 32.7311031794905
$ws.Range("L6").Value = 0.0663879265387853
$ws.Range("M6").Value = 0.201015738149484

$ws.Range("A7").Value = "Derivation"
$ws.Range("B7").Value = 0.646708333333333
$ws.Range("C7").Value = 0.6395
$ws.Range("D7").Formula = "=ABS(B7-C7)"
$ws.Range("E7").Value = 34.6644291666666
$ws.Range("F7").Value = 34.6640666666666
$ws.Range("G7").Formula = "=ABS(E7-F7)"
$ws.Range("H7").Value = 0.9837461031245034
$ws.Range("I7").Value = 0.983745416967042
$ws.Range("J7").Formula = "=ABS(H7-I7)"
$ws.Range("K7").Value = 31.0610796468515
$ws.Range("L7").Value = 0.0603809903065363
$ws.Range("M7").Value = 0.194030654927094

$ws.Range("A8").Value = "Derivation"
$ws.Range("B8").Value = 0.71475
$ws.Range("C8").Value = 0.706666666666666
$ws.Range("D8").Formula = "=ABS(B8-C8)"
$ws.Range("E8").Value = 34.8418625
$ws.Range("F8").Value = 34.8414916666666
$ws.Range("G8").Formula = "=ABS(E8-F8)"
$ws.Range("H8").Value = 0.9844040260899775
$ws.Range("I8").Value = 0.9844020059621982
$ws.Range("J8").Formula = "=ABS(H8-I8)"
$ws.Range("K8").Value = 30.2891943100195
$ws.Range("L8").Value = 0.0588848143815994
$ws.Range("M8").Value = 0.19218107809623

# The extra "Derivation"/"Anchor" column-A text is wider than the old lambda values,
# and column H got a touch wider too -- mirror the resulting manual column resize
# (this also drops the bestFit flag on those two columns, same as the real edit).
$ws.Columns.Item(1).ColumnWidth = 9.36
$ws.Columns.Item(8).ColumnWidth = 9.8

# Match the saved print setup or the new sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the same selection/active-cell state the author had when saving.
$ws.Range("J17").Select()
$ws.Activate()
